# Refresh the cryptos price/volume snapshot (GitHub Actions-style data pull).
# Column D ("Price") and column E ("Volume(1h)") are plain text cells in the
# source sheet (values such as "30.445.57" or "  +1.56%  " are not numbers),
# so cells whose new text would otherwise auto-parse as a number are written
# with a leading apostrophe (Excel's text quote-prefix) to keep them as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.445.57'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '1.997.19'
$ws.Range('E3').Value = '  +4.01%  '
$ws.Range('D5').Value = "'324.61"
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = "'0.5110"
$ws.Range('E7').Value = '  +1.47%  '
$ws.Range('D8').Value = "'0.4135"
$ws.Range('E8').Value = '  +2.59%  '
$ws.Range('D9').Value = "'0.08716"
$ws.Range('E9').Value = '  +5.79%  '
$ws.Range('D10').Value = "'1.131"
$ws.Range('D11').Value = "'43.04"
$ws.Range('E11').Value = '  +2.40%  '
$ws.Range('D12').Value = "'24.67"
$ws.Range('E12').Value = '  +4.48%  '
$ws.Range('D13').Value = '1.991.95'
$ws.Range('E13').Value = '  +4.11%  '
$ws.Range('D14').Value = "'6.566"
$ws.Range('E14').Value = '  +2.54%  '
$ws.Range('D15').Value = "'7.433"
$ws.Range('E15').Value = '  +1.93%  '
$ws.Range('D16').Value = "'1.002"
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').Value = "'94.39"
$ws.Range('E17').Value = '  +2.45%  '
$ws.Range('D18').Value = "'0.00001116"
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('E20').Value = '  +3.71%  '
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').Value = "'6.166"
$ws.Range('E22').Value = '  +3.56%  '
$ws.Range('D23').Value = '30.494.43'
$ws.Range('E23').Value = '  +1.58%  '
$ws.Range('E24').Value = '  +4.74%  '
$ws.Range('D25').Value = "'2.226"
$ws.Range('E25').Value = '  +1.47%  '
$ws.Range('D26').Value = '2.228.89'
$ws.Range('E26').Value = '  +4.47%  '
$ws.Range('D27').Value = "'22.38"
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('D28').Value = "'163.42"
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('D29').Value = "'2.398"
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('D30').Value = "'131.51"
$ws.Range('E30').Value = '  +1.98%  '
$ws.Range('D31').Value = "'1.136"
$ws.Range('E31').Value = '  +0.76%  '
$ws.Range('D32').Value = "'0.1052"
$ws.Range('E32').Value = '  +0.96%  '
$ws.Range('D33').Value = "'6.073"
$ws.Range('E33').Value = '  +1.52%  '
$ws.Range('D34').Value = "'3.854"
$ws.Range('E34').Value = '  +1.01%  '
$ws.Range('E35').Value = '  +12.13%  '
$ws.Range('D36').Value = "'0.02520"
$ws.Range('E36').Value = '  +3.06%  '
$ws.Range('D37').Value = "'5.417"
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('D38').Value = "'0.06608"
$ws.Range('E38').Value = '  +2.52%  '
$ws.Range('D39').Value = "'12.35"
$ws.Range('E39').Value = '  +8.45%  '
$ws.Range('E40').Value = '  +1.86%  '
$ws.Range('D41').Value = "'8.990"
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').Value = "'0.6621"
$ws.Range('E42').Value = '  +3.09%  '
$ws.Range('D43').Value = "'1.235"
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('D44').Value = "'13.63"
$ws.Range('E44').Value = '  +2.41%  '
$ws.Range('D45').Value = "'0.6158"
$ws.Range('E45').Value = '  +2.59%  '
$ws.Range('D46').Value = "'2.201"
$ws.Range('D47').Value = "'3.667"
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('D48').Value = "'1.265"
$ws.Range('E48').Value = '  +4.08%  '
$ws.Range('D49').Value = "'124.42"
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('D50').Value = "'80.11"
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('D51').Value = "'0.06893"
$ws.Range('E51').Value = '  +1.57%  '
